$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("F12").Value = "Nhi, Tú"
$v2 = $ws.Cells.Item(12,6).Value2
Write-Output "F12 value2: $v2"
